# Populate the team-specific Markov transition matrix on Sheet1.
# Rows/columns are game states (Af0..Af3, Ai0..Ai3, Ar0, Bf0..Bf3, Bi0..Bi3, Br0);
# each row holds the observed transition-probability distribution out of that
# state now that more simulated games have been folded into the counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Af0)
$ws.Range("B2").Value = 0.3636363636363636
$ws.Range("C2").Value = 0.2727272727272727
$ws.Range("P2").Value = 0.2727272727272727
$ws.Range("S2").Value = 0.09090909090909091

# Row 3 (Af1)
$ws.Range("P3").Value = 0.6666666666666666
$ws.Range("S3").Value = 0.3333333333333333

# Row 6 (Ai0)
$ws.Range("B6").Value = 0.1176470588235294
$ws.Range("F6").Value = 0.1176470588235294
$ws.Range("J6").Value = 0.2941176470588235
$ws.Range("Q6").Value = 0.2941176470588235
$ws.Range("R6").Value = 0.05882352941176471
$ws.Range("S6").Value = 0.1176470588235294

# Row 8 (Ai2)
$ws.Range("B8").Value = 0.09090909090909091
$ws.Range("F8").Value = 0.1818181818181818
$ws.Range("J8").Value = 0.09090909090909091
$ws.Range("Q8").Value = 0.1818181818181818
$ws.Range("R8").Value = 0.1818181818181818
$ws.Range("S8").Value = 0.2727272727272727

# Row 9 (Ai3)
$ws.Range("F9").Value = 0.1428571428571428
$ws.Range("O9").Value = 0.1428571428571428
$ws.Range("R9").Value = 0.1428571428571428
$ws.Range("S9").Value = 0.5714285714285714

# Row 10 (Ar0)
$ws.Range("B10").Value = 0.06557377049180328
$ws.Range("F10").Value = 0.06557377049180328
$ws.Range("J10").Value = 0.1147540983606557
$ws.Range("O10").Value = 0.03278688524590164
$ws.Range("Q10").Value = 0.3770491803278688
$ws.Range("R10").Value = 0.09836065573770492
$ws.Range("S10").Value = 0.2459016393442623

# Row 11 (Bf0)
$ws.Range("J11").Value = 0.5
$ws.Range("L11").Value = 0.5

# Row 12 (Bf1)
$ws.Range("J12").Value = 1

# Row 15 (Bi0)
$ws.Range("H15").Value = 0.1111111111111111
$ws.Range("I15").Value = 0.2222222222222222
$ws.Range("J15").Value = 0.3333333333333333
$ws.Range("O15").Value = 0.1111111111111111
$ws.Range("S15").Value = 0.2222222222222222

# Row 16 (Bi1)
$ws.Range("F16").Value = 0.4
$ws.Range("J16").Value = 0.6

# Row 17 (Bi2)
$ws.Range("F17").Value = 0.03333333333333333
$ws.Range("H17").Value = 0.1333333333333333
$ws.Range("I17").Value = 0.1
$ws.Range("J17").Value = 0.6333333333333333
$ws.Range("K17").Value = 0.03333333333333333
$ws.Range("O17").Value = 0.03333333333333333
$ws.Range("S17").Value = 0.03333333333333333

# Row 18 (Bi3)
$ws.Range("F18").Value = 0.2
$ws.Range("H18").Value = 0.1
$ws.Range("I18").Value = 0.1
$ws.Range("J18").Value = 0.6

# Row 19 (Br0)
$ws.Range("F19").Value = 0.06896551724137931
$ws.Range("H19").Value = 0.2068965517241379
$ws.Range("I19").Value = 0.03448275862068965
$ws.Range("J19").Value = 0.5172413793103449
$ws.Range("K19").Value = 0.03448275862068965
$ws.Range("O19").Value = 0.1379310344827586
